$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 114-141 (Remessa, Material, Quantidade).
# Column A ("Remessa") values are pure numeric-looking IDs that must be
# stored as TEXT (matching the rest of the column), so they are written
# through a scratch cell that is temporarily formatted as Text ("@") and
# then copied across as values - this avoids Excel's automatic
# string->number coercion while keeping the destination cell's own
# pre-existing style/formatting untouched.
$data = @(
    @(114, "80266408", "10386-ARI-I", 1),
    @(115, "80266409", "20087-CTY-I", 1),
    @(116, "80266410", "10493-ARI-I", 1),
    @(117, "80266411", "10493-ARI-I", 1),
    @(118, "80266412", "21481-NZX-I", 1),
    @(119, "80266413", "10255-ARI-I", 4),
    @(120, "80266413", "10258-ARI-I", 4),
    @(121, "80266413", "10259-ARI-I", 1),
    @(122, "80266413", "10256-ARI-I", 4),
    @(123, "80266413", "10257-ARI-I", 1),
    @(124, "80266413", "10493-ARI-I", 8),
    @(125, "80266413", "10636-ARI-I", 4),
    @(126, "80266413", "10637-ARI-I", 2),
    @(127, "80266413", "10486-ARI-I", 1),
    @(128, "80266413", "10489-ARI-I", 1),
    @(129, "80266413", "10487-ARI-I", 2),
    @(130, "80266413", "10479-ARI-I", 1),
    @(131, "80266413", "10246-ARI-I", 1),
    @(132, "80266413", "10195-ARI-I", 10),
    @(133, "80266414", "21340-NZX-I", 1),
    @(134, "80266418", "10369-ARI-I", 1),
    @(135, "80266419", "14099-TDK-N", 1000),
    @(136, "80266424", "11848-KRO-L", 1),
    @(137, "80266425", "10247-ARI-I", 1),
    @(138, "80266426", "21475-NZX-I", 1),
    @(139, "80266426", "23359-ATE-I", 2),
    @(140, "80266427", "11820-DLO-I", 1),
    @(141, "80266427", "11848-DLO-I", 2)
)

# Scratch cell, far outside the used range, used to mint text-typed values.
$scratch = $ws.Cells.Item(1, 20)
$scratch.NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]

    # Column A: force text storage via the scratch cell + paste-values trick.
    $scratch.Value = $row[1]
    $scratch.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)  # xlPasteValues

    # Column B: already text because it contains letters, plain assignment
    # keeps it as text and keeps the destination's existing style.
    $ws.Cells.Item($r, 2).Value = $row[2]

    # Column C: plain numeric quantity.
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Clean up the scratch cell/clipboard state.
$scratch.Clear()
$excel.CutCopyMode = $false

# Update the sheet view: selection now spans A2:C141, and the previous
# topLeftCell="A82" scroll position is reset (sheet view scrolled back to top).
$ws.Range("A2:C141").Select()
$excel.ActiveWindow.ScrollRow = 1

$wb.Save()
